# Update column C (Num_Inclusions) values in Sheet1 to use 3rd quartile instead of mean
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 2
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 12
$ws.Range("C10").Value = 5
$ws.Range("C12").Value = 11
$ws.Range("C16").Value = 3
$ws.Range("C17").Value = 0
$ws.Range("C19").Value = 9
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 14
